# Add a new input variable "MDTOUT" (outlet mass flow rate) as a new row
# in the CHAN sheet of template_conductor_1_operation.xlsx, right below the
# existing MDTIN (inlet mass flow rate) row, and make CHAN the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHAN")

# Insert a new row at position 12 (pushes the existing FLOWDIR row, and
# everything below it, down by one row) and fill it in with the new
# MDTOUT variable definition, matching the layout of the MDTIN row above it.
$ws.Rows.Item(12).Insert()

$ws.Range("A12").Value = "MDTOUT"
$ws.Range("B12").Value = "kg/s"
$ws.Range("C12").Value = "float"
$ws.Range("D12").Value = "outlet mass flow rate "
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0

# Make CHAN the active sheet/tab and leave the active cell at G12, matching
# the cursor position after entering the new row of data.
$ws.Activate()
$ws.Range("G12").Select()
